# Update the cryptos price (column D) and volume/1h change (column E) values
# on the active worksheet, per the latest GitHub Actions scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '27.020.84'; E = '  -1.28%  ' },
    @{ Row = 3; D = '1.824.94'; E = '  -0.17%  ' },
    @{ Row = 4; D = '1.001'; E = '  -0.22%  ' },
    @{ Row = 5; D = '311.52'; E = '  -1.28%  ' },
    @{ Row = 6; D = $null; E = '  -0.30%  ' },
    @{ Row = 7; D = '0.4400'; E = '  +2.34%  ' },
    @{ Row = 8; D = '0.3675'; E = '  -0.70%  ' },
    @{ Row = 9; D = '0.07273'; E = '  +0.22%  ' },
    @{ Row = 10; D = '0.8432'; E = '  -2.65%  ' },
    @{ Row = 11; D = '20.68'; E = '  -2.27%  ' },
    @{ Row = 12; D = '1.812.29'; E = '  -0.82%  ' },
    @{ Row = 13; D = '6.646'; E = '  -0.51%  ' },
    @{ Row = 14; D = '0.07075'; E = '  -0.40%  ' },
    @{ Row = 15; D = '5.299'; E = '  -1.06%  ' },
    @{ Row = 16; D = '89.56'; E = '  +1.95%  ' },
    @{ Row = 17; D = '1.003'; E = '  -0.30%  ' },
    @{ Row = 18; D = '0.000008779'; E = '  -1.32%  ' },
    @{ Row = 19; D = $null; E = '  -0.21%  ' },
    @{ Row = 20; D = $null; E = '  -2.04%  ' },
    @{ Row = 21; D = '26.812.66'; E = '  -2.08%  ' },
    @{ Row = 22; D = '5.146'; E = '  -0.31%  ' },
    @{ Row = 23; D = '10.89'; E = '  +0.20%  ' },
    @{ Row = 24; D = '2.050.50'; E = '  -0.05%  ' },
    @{ Row = 25; D = '1.981'; E = '  -1.53%  ' },
    @{ Row = 26; D = '151.58'; E = '  -1.10%  ' },
    @{ Row = 27; D = '2.207'; E = '  +3.13%  ' },
    @{ Row = 28; D = '18.27'; E = '  -1.08%  ' },
    @{ Row = 29; D = '5.227'; E = '  -1.35%  ' },
    @{ Row = 30; D = '117.09'; E = '  +0.10%  ' },
    @{ Row = 31; D = '0.08780'; E = '  -0.88%  ' },
    @{ Row = 32; D = $null; E = '  -2.39%  ' },
    @{ Row = 33; D = '0.7396'; E = '  -3.51%  ' },
    @{ Row = 34; D = $null; E = '  +2.03%  ' },
    @{ Row = 35; D = '4.424'; E = '  -1.83%  ' },
    @{ Row = 36; D = '1.0000'; E = '  -0.39%  ' },
    @{ Row = 37; D = $null; E = '  -2.20%  ' },
    @{ Row = 38; D = '0.01947'; E = '  -0.66%  ' },
    @{ Row = 39; D = '0.05241'; E = '  -0.60%  ' },
    @{ Row = 40; D = '7.243'; E = '  +1.60%  ' },
    @{ Row = 41; D = '2.868'; E = '  -0.45%  ' },
    @{ Row = 42; D = '0.5150'; E = '  +1.62%  ' },
    @{ Row = 43; D = '0.1695'; E = '  +0.79%  ' },
    @{ Row = 44; D = $null; E = '  -1.58%  ' },
    @{ Row = 45; D = '10.57'; E = '  +0.05%  ' },
    @{ Row = 46; D = '0.4809'; E = '  +1.65%  ' },
    @{ Row = 47; D = '106.03'; E = '  -0.30%  ' },
    @{ Row = 48; D = '1.926'; E = '  +5.50%  ' },
    @{ Row = 49; D = '0.9999'; E = '  -0.37%  ' },
    @{ Row = 50; D = '0.06334'; E = '  -1.34%  ' },
    @{ Row = 51; D = '1.657'; E = '  -0.96%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
